$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 58, pushing the existing row 58 (and everything below
# it, through the former row 153) down by one -- the former row 153 becomes
# the new row 154, matching the dimension change to A1:T154.
$ws.Rows.Item(58).Insert()

# Populate the newly-inserted row 58 with the new price-report record.
$ws.Cells.Item(58, 1).Value = 2
$ws.Cells.Item(58, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(58, 3).Value = "Coquimbo"
$ws.Cells.Item(58, 4).Value = 44665
$ws.Cells.Item(58, 5).Value = 4
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100109
$ws.Cells.Item(58, 8).Value = "Uva"
$ws.Cells.Item(58, 9).Value = 100109001
$ws.Cells.Item(58, 10).Value = "Uva"
$ws.Cells.Item(58, 11).Value = "Red Globe"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 400
$ws.Cells.Item(58, 14).Value = 8000
$ws.Cells.Item(58, 15).Value = 8500
$ws.Cells.Item(58, 16).Value = 8250
$ws.Cells.Item(58, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(58, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(58, 19).Value = 458
$ws.Cells.Item(58, 20).Value = 18
